# "Implemented Label wise bill generate(Individual)"
# Fill in the bill-header labels with this teacher's actual details, enter
# the quantities that were previously left blank, and write the amount
# in words once the total is known.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header labels: append the individual's details after the fixed label text ---
$ws.Range("A3").Value = "নাম: Dr. Sk. Imran Hossain"
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# --- Quantities for this individual's bill (previously blank) ---
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("G28").Value = 1

# --- Amount in words, once the grand total (I32) is known ---
$ws.Range("A32").Value = "কথায়:তের হাজার ছয়শত পঁচান্ন টাকা মাত্র।"

# --- Leave the selection where the author left it while reviewing the sheet ---
$null = $ws.Range("B5").Select()
